# Configuration update after IOS 2015
#
# This script reproduces, via Excel COM-interop calls, the edits captured
# in the target OOXML diff:
#   1. "icf.185"   (sheet 1) - selection moved from BS2:BU2 to BT27
#   2. "icf_a.185" (sheet 2, the active tab) - view scrolled/selection moved
#      to BU12, the BO1:BU1 date-header formulas re-entered as one fill
#      (producing a shared formula group), and two data points updated
#      (BU9, BU12).
#   3. "Eventos.185" (sheet 3) - the A66:A71 date formulas re-entered as one
#      fill (producing a shared formula group) for the rows added after the
#      IOS 2015 update.

$wb = $excel.ActiveWorkbook

$wsIcf   = $wb.Worksheets.Item(1)   # icf.185
$wsIcfA  = $wb.Worksheets.Item(2)   # icf_a.185 (active tab)
$wsEvt   = $wb.Worksheets.Item(3)   # Eventos.185

# --- 1. icf.185: move the selection (sheet is not the active tab, so we
#        briefly activate it, select, then restore the original active
#        sheet at the end). -------------------------------------------------
$wsIcf.Activate()
$wsIcf.Range("BT27").Select()

# --- 2. icf_a.185: rebuild the BO1:BU1 formulas as a single fill operation
#        so they collapse into a shared formula group (matches the diff,
#        same underlying formula / same computed values), then update the
#        two changed data cells, then move the selection. ------------------
$wsIcfA.Activate()
$wsIcfA.Range("BO1:BU1").Formula = "=+BO2-693960"

$wsIcfA.Range("BU9").Value = 0.342
$wsIcfA.Range("BU12").Value = 1570

$wsIcfA.Range("BU12").Select()

# --- 3. Eventos.185: rebuild the A66:A71 formulas as a single fill
#        operation so they collapse into a shared formula group (matches
#        the diff, same underlying formula / same computed values). -------
$wsEvt.Range("A66:A71").Formula = "=+B66-693960"

# Restore icf_a.185 as the active sheet/tab (it was active before this
# script ran, and remains so afterwards per the diff).
$wsIcfA.Activate()
